$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 11.75984487686059
$ws.Cells.Item(2, 3).Value = 11.25766534996142
$ws.Cells.Item(2, 4).Value = 6.003545107623809
$ws.Cells.Item(2, 5).Value = 13.43050152248331
$ws.Cells.Item(2, 7).Value = 36.61453259059855
$ws.Cells.Item(2, 8).Value = 16.33255314349849
$ws.Cells.Item(2, 11).Value = 8.163292196487014
$ws.Cells.Item(2, 12).Value = 9.944229855399318
$ws.Cells.Item(2, 13).Value = 13.81078491062013
$ws.Cells.Item(2, 15).Value = 25.8940832894901
$ws.Cells.Item(3, 2).Value = 11.53013456289844
$ws.Cells.Item(3, 3).Value = 11.26755525584063
$ws.Cells.Item(3, 4).Value = 5.887808091363695
$ws.Cells.Item(3, 5).Value = 13.46358376077121
$ws.Cells.Item(3, 7).Value = 36.69210125673893
$ws.Cells.Item(3, 8).Value = 16.38002733073477
$ws.Cells.Item(3, 11).Value = 7.980981521310334
$ws.Cells.Item(3, 12).Value = 9.951702718769191
$ws.Cells.Item(3, 13).Value = 13.77837831260453
$ws.Cells.Item(3, 15).Value = 25.96953070275308
$ws.Cells.Item(4, 2).Value = 11.38884136864393
$ws.Cells.Item(4, 3).Value = 11.27417789836512
$ws.Cells.Item(4, 4).Value = 5.817322312425629
$ws.Cells.Item(4, 5).Value = 13.48574196875131
$ws.Cells.Item(4, 7).Value = 36.74940450299063
$ws.Cells.Item(4, 8).Value = 16.41157924747806
$ws.Cells.Item(4, 11).Value = 7.868044384719405
$ws.Cells.Item(4, 12).Value = 9.957612756420435
$ws.Cells.Item(4, 13).Value = 13.76042088930586
$ws.Cells.Item(4, 15).Value = 26.02077372946392
$ws.Cells.Item(5, 2).Value = 11.33128010793875
$ws.Cells.Item(5, 3).Value = 11.27701547338697
$ws.Cells.Item(5, 4).Value = 5.788788126726148
$ws.Cells.Item(5, 5).Value = 13.49523592726935
$ws.Cells.Item(5, 7).Value = 36.77518212331665
$ws.Cells.Item(5, 8).Value = 16.42504110842975
$ws.Cells.Item(5, 11).Value = 7.821835246801857
$ws.Cells.Item(5, 12).Value = 9.960354030462163
$ws.Cells.Item(5, 13).Value = 13.75359627839503
$ws.Cells.Item(5, 15).Value = 26.04289066905024
$ws.Cells.Item(6, 2).Value = 11.32172548606094
$ws.Cells.Item(6, 3).Value = 11.27749504641391
$ws.Cells.Item(6, 4).Value = 5.784062780616382
$ws.Cells.Item(6, 5).Value = 13.49684044129009
$ws.Cells.Item(6, 7).Value = 36.7796087880248
$ws.Cells.Item(6, 8).Value = 16.42731293379129
$ws.Cells.Item(6, 11).Value = 7.814152932158462
$ws.Cells.Item(6, 12).Value = 9.960829336984659
$ws.Cells.Item(6, 13).Value = 13.75249299083396
$ws.Cells.Item(6, 15).Value = 26.04663770953776
$ws.Cells.Item(7, 2).Value = 11.38806490246525
$ws.Cells.Item(7, 3).Value = 11.27421560433277
$ws.Cells.Item(7, 4).Value = 5.816936665128004
$ws.Cells.Item(7, 5).Value = 13.48586812733572
$ws.Cells.Item(7, 7).Value = 36.74974233568998
$ws.Cells.Item(7, 8).Value = 16.41175835240531
$ws.Cells.Item(7, 11).Value = 7.867421859599949
$ws.Cells.Item(7, 12).Value = 9.957648377688585
$ws.Cells.Item(7, 13).Value = 13.76032684664346
$ws.Cells.Item(7, 15).Value = 26.02106700811865
$ws.Cells.Item(8, 2).Value = 11.68073997079394
$ws.Cells.Item(8, 3).Value = 11.26096145994407
$ws.Cells.Item(8, 4).Value = 5.963545564203976
$ws.Cells.Item(8, 5).Value = 13.44152536585524
$ws.Cells.Item(8, 7).Value = 36.63926609865904
$ws.Cells.Item(8, 8).Value = 16.34842362281221
$ws.Cells.Item(8, 11).Value = 8.100676487973546
$ws.Cells.Item(8, 12).Value = 9.946532528010888
$ws.Cells.Item(8, 13).Value = 13.79921165794972
$ws.Cells.Item(8, 15).Value = 25.9190758379073
$ws.Cells.Item(9, 2).Value = 12.24928770342155
$ws.Cells.Item(9, 3).Value = 11.239316399204
$ws.Cells.Item(9, 4).Value = 6.253712474446947
$ws.Cells.Item(9, 5).Value = 13.36920309100168
$ws.Cells.Item(9, 7).Value = 36.49967452783036
$ws.Cells.Item(9, 8).Value = 16.24328530793875
$ws.Cells.Item(9, 11).Value = 8.547469410355896
$ws.Cells.Item(9, 12).Value = 9.935195797227211
$ws.Cells.Item(9, 13).Value = 13.8906149423153
$ws.Cells.Item(9, 15).Value = 25.75816627420974
$ws.Cells.Item(10, 2).Value = 12.65934234957739
$ws.Cells.Item(10, 3).Value = 11.22603755718194
$ws.Cells.Item(10, 4).Value = 6.466033415426374
$ws.Cells.Item(10, 5).Value = 13.32497494498386
$ws.Cells.Item(10, 7).Value = 36.44442202024636
$ws.Cells.Item(10, 8).Value = 16.17765664783385
$ws.Cells.Item(10, 11).Value = 8.865846460581373
$ws.Cells.Item(10, 12).Value = 9.9332075928361
$ws.Cells.Item(10, 13).Value = 13.96665558145867
$ws.Cells.Item(10, 15).Value = 25.6638744765734
$ws.Cells.Item(11, 2).Value = 12.84336162808595
$ws.Cells.Item(11, 3).Value = 11.22056097546352
$ws.Cells.Item(11, 4).Value = 6.561937210918146
$ws.Cells.Item(11, 5).Value = 13.30678528180057
$ws.Cells.Item(11, 7).Value = 36.42961088292196
$ws.Cells.Item(11, 8).Value = 16.15032210250898
$ws.Cells.Item(11, 11).Value = 9.007883942235797
$ws.Cells.Item(11, 12).Value = 9.933670812102156
$ws.Cells.Item(11, 13).Value = 14.00309668214421
$ws.Cells.Item(11, 15).Value = 25.62619304317847
$ws.Cells.Item(12, 2).Value = 12.91261612541975
$ws.Cells.Item(12, 3).Value = 11.21856780552615
$ws.Cells.Item(12, 4).Value = 6.598116345797635
$ws.Cells.Item(12, 5).Value = 13.30017458457236
$ws.Cells.Item(12, 7).Value = 36.42548947872058
$ws.Cells.Item(12, 8).Value = 16.14033360230207
$ws.Cells.Item(12, 11).Value = 9.061218323828852
$ws.Cells.Item(12, 12).Value = 9.934041974371976
$ws.Cells.Item(12, 13).Value = 14.01715455088374
$ws.Cells.Item(12, 15).Value = 25.61267495952634
$ws.Cells.Item(13, 2).Value = 12.89772107292905
$ws.Cells.Item(13, 3).Value = 11.21899348807164
$ws.Cells.Item(13, 4).Value = 6.590331245382942
$ws.Cells.Item(13, 5).Value = 13.30158598522643
$ws.Cells.Item(13, 7).Value = 36.42631091802006
$ws.Cells.Item(13, 8).Value = 16.14246868097375
$ws.Cells.Item(13, 11).Value = 9.049752682248558
$ws.Cells.Item(13, 12).Value = 9.933953347257027
$ws.Cells.Item(13, 13).Value = 14.01411556829536
$ws.Cells.Item(13, 15).Value = 25.61555289427508
$ws.Cells.Item(14, 2).Value = 12.84906822093311
$ws.Cells.Item(14, 3).Value = 11.22039538133356
$ws.Cells.Item(14, 4).Value = 6.564916665430053
$ws.Cells.Item(14, 5).Value = 13.30623585873264
$ws.Cells.Item(14, 7).Value = 36.42924199365819
$ws.Cells.Item(14, 8).Value = 16.14949307708846
$ws.Cells.Item(14, 11).Value = 9.012281130798115
$ws.Cells.Item(14, 12).Value = 9.933697430221146
$ws.Cells.Item(14, 13).Value = 14.00424809698942
$ws.Cells.Item(14, 15).Value = 25.62506584170134
$ws.Cells.Item(15, 2).Value = 12.81920905483675
$ws.Cells.Item(15, 3).Value = 11.22126457776186
$ws.Cells.Item(15, 4).Value = 6.549330452647768
$ws.Cells.Item(15, 5).Value = 13.3091201515979
$ws.Cells.Item(15, 7).Value = 36.43123111039716
$ws.Cells.Item(15, 8).Value = 16.15384293308544
$ws.Cells.Item(15, 11).Value = 8.98926841029359
$ws.Cells.Item(15, 12).Value = 9.933566137724448
$ws.Cells.Item(15, 13).Value = 13.99823740873758
$ws.Cells.Item(15, 15).Value = 25.63099064965311
$ws.Cells.Item(16, 2).Value = 12.64725973599892
$ws.Cells.Item(16, 3).Value = 11.22640677894092
$ws.Cells.Item(16, 4).Value = 6.459748730263826
$ws.Cells.Item(16, 5).Value = 13.32620250177164
$ws.Cells.Item(16, 7).Value = 36.44559795589913
$ws.Cells.Item(16, 8).Value = 16.17949374647215
$ws.Cells.Item(16, 11).Value = 8.856503367672266
$ws.Cells.Item(16, 12).Value = 9.933204764573254
$ws.Cells.Item(16, 13).Value = 13.96431065145031
$ws.Cells.Item(16, 15).Value = 25.66644206627086
$ws.Cells.Item(17, 2).Value = 12.54108038674736
$ws.Cells.Item(17, 3).Value = 11.22970551893088
$ws.Cells.Item(17, 4).Value = 6.404589811904566
$ws.Cells.Item(17, 5).Value = 13.33717614496787
$ws.Cells.Item(17, 7).Value = 36.45705792020408
$ws.Cells.Item(17, 8).Value = 16.19587525414943
$ws.Cells.Item(17, 11).Value = 8.774303737523631
$ws.Cells.Item(17, 12).Value = 9.933332736620413
$ws.Cells.Item(17, 13).Value = 13.94396616956331
$ws.Cells.Item(17, 15).Value = 25.6895264997917
$ws.Cells.Item(18, 2).Value = 12.4797757588071
$ws.Cells.Item(18, 3).Value = 11.23165598069973
$ws.Cells.Item(18, 4).Value = 6.372801895688156
$ws.Cells.Item(18, 5).Value = 13.34366956199095
$ws.Cells.Item(18, 7).Value = 36.46462098482562
$ws.Cells.Item(18, 8).Value = 16.2055346781224
$ws.Cells.Item(18, 11).Value = 8.726764495838349
$ws.Cells.Item(18, 12).Value = 9.93353513334608
$ws.Cells.Item(18, 13).Value = 13.93243915076289
$ws.Cells.Item(18, 15).Value = 25.70329458546221
$ws.Cells.Item(19, 2).Value = 12.45898116285258
$ws.Cells.Item(19, 3).Value = 11.23232550854564
$ws.Cells.Item(19, 4).Value = 6.362029647539716
$ws.Cells.Item(19, 5).Value = 13.34589932853806
$ws.Cells.Item(19, 7).Value = 36.467348481045
$ws.Cells.Item(19, 8).Value = 16.20884593205209
$ws.Cells.Item(19, 11).Value = 8.710625432946371
$ws.Cells.Item(19, 12).Value = 9.933625808017657
$ws.Cells.Item(19, 13).Value = 13.92856651230464
$ws.Cells.Item(19, 15).Value = 25.70804042332421
$ws.Cells.Item(20, 2).Value = 12.5524079852841
$ws.Cells.Item(20, 3).Value = 11.22934886836047
$ws.Cells.Item(20, 4).Value = 6.410468278667603
$ws.Cells.Item(20, 5).Value = 13.33598918131653
$ws.Cells.Item(20, 7).Value = 36.45573740776047
$ws.Cells.Item(20, 8).Value = 16.19410686325777
$ws.Cells.Item(20, 11).Value = 8.783081333589337
$ws.Cells.Item(20, 12).Value = 9.933305790841603
$ws.Cells.Item(20, 13).Value = 13.94611385907984
$ws.Cells.Item(20, 15).Value = 25.68701834297437
$ws.Cells.Item(21, 2).Value = 12.86337092224328
$ws.Cells.Item(21, 3).Value = 11.21998142461086
$ws.Cells.Item(21, 4).Value = 6.572385580818585
$ws.Cells.Item(21, 5).Value = 13.30486255301631
$ws.Cells.Item(21, 7).Value = 36.4283406852537
$ws.Cells.Item(21, 8).Value = 16.14742000289181
$ws.Cells.Item(21, 11).Value = 9.023300087775796
$ws.Cells.Item(21, 12).Value = 9.933767294310366
$ws.Cells.Item(21, 13).Value = 14.00713946117532
$ws.Cells.Item(21, 15).Value = 25.62225126342727
$ws.Cells.Item(22, 2).Value = 13.06407052257598
$ws.Cells.Item(22, 3).Value = 11.21432940713268
$ws.Cells.Item(22, 4).Value = 6.677389351076005
$ws.Cells.Item(22, 5).Value = 13.28613591350909
$ws.Cells.Item(22, 7).Value = 36.41910473875808
$ws.Cells.Item(22, 8).Value = 16.11902039731285
$ws.Cells.Item(22, 11).Value = 9.177640786816896
$ws.Cells.Item(22, 12).Value = 9.93520942245981
$ws.Cells.Item(22, 13).Value = 14.0485257300272
$ws.Cells.Item(22, 15).Value = 25.58430058925865
$ws.Cells.Item(23, 2).Value = 12.95720643891705
$ws.Cells.Item(23, 3).Value = 11.21730311275257
$ws.Cells.Item(23, 4).Value = 6.621434248802038
$ws.Cells.Item(23, 5).Value = 13.29598283880749
$ws.Cells.Item(23, 7).Value = 36.42324024346298
$ws.Cells.Item(23, 8).Value = 16.13398443835685
$ws.Cells.Item(23, 11).Value = 9.09552498496029
$ws.Cells.Item(23, 12).Value = 9.934335697372438
$ws.Cells.Item(23, 13).Value = 14.02630220379002
$ws.Cells.Item(23, 15).Value = 25.60415448695599
$ws.Cells.Item(24, 2).Value = 12.54728758483349
$ws.Cells.Item(24, 3).Value = 11.22950994193423
$ws.Cells.Item(24, 4).Value = 6.407810858549054
$ws.Cells.Item(24, 5).Value = 13.33652523280074
$ws.Cells.Item(24, 7).Value = 36.45633137589118
$ws.Cells.Item(24, 8).Value = 16.19490560058899
$ws.Cells.Item(24, 11).Value = 8.779113855598284
$ws.Cells.Item(24, 12).Value = 9.933317571689741
$ws.Cells.Item(24, 13).Value = 13.94514236074252
$ws.Cells.Item(24, 15).Value = 25.68815073408915
$ws.Cells.Item(25, 2).Value = 12.09652057481058
$ws.Cells.Item(25, 3).Value = 11.24470934242761
$ws.Cells.Item(25, 4).Value = 6.175191431766942
$ws.Cells.Item(25, 5).Value = 13.38720294245529
$ws.Cells.Item(25, 7).Value = 36.52915094179532
$ws.Cells.Item(25, 8).Value = 16.26968776122887
$ws.Cells.Item(25, 11).Value = 8.428101451402945
$ws.Cells.Item(25, 12).Value = 9.937146221302276
$ws.Cells.Item(25, 13).Value = 13.86430066223908
$ws.Cells.Item(25, 15).Value = 25.79750097638654
